# Update the 25 "two-digit / one-digit" division problems in the single
# table on the page. Each data row is followed by 3 blank rows, and the
# mapping of old -> new values is positional (several old values repeat,
# e.g. "84÷7=", and several new values repeat too, e.g. "70÷6=" and
# "80÷3="), so we must address each cell explicitly by (row, column)
# rather than doing a global text replace.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-Problem {
    param($Row, $Col, $Old, $New)

    $cell = $t.Cell($Row, $Col)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $Old) {
        throw "Cell ($Row,$Col) expected '$Old' but found '$current'"
    }
    $cell.Range.Text = $New
}

# Row 1 (table row 1)
Set-Problem 1 1 "46÷7=" "13÷2="
Set-Problem 1 2 "30÷8=" "36÷7="
Set-Problem 1 3 "45÷9=" "91÷5="
Set-Problem 1 4 "17÷8=" "77÷4="
Set-Problem 1 5 "89÷5=" "33÷8="

# Row 2 (table row 5)
Set-Problem 5 1 "44÷5=" "21÷5="
Set-Problem 5 2 "79÷9=" "19÷3="
Set-Problem 5 3 "26÷4=" "56÷9="
Set-Problem 5 4 "98÷5=" "78÷7="
Set-Problem 5 5 "47÷9=" "73÷4="

# Row 3 (table row 9)
Set-Problem 9 1 "14÷4=" "80÷3="
Set-Problem 9 2 "84÷7=" "75÷8="
Set-Problem 9 3 "57÷7=" "78÷7="
Set-Problem 9 4 "73÷6=" "88÷5="
Set-Problem 9 5 "34÷6=" "52÷8="

# Row 4 (table row 13)
Set-Problem 13 1 "72÷7=" "91÷3="
Set-Problem 13 2 "85÷5=" "70÷6="
Set-Problem 13 3 "44÷2=" "44÷4="
Set-Problem 13 4 "89÷7=" "70÷6="
Set-Problem 13 5 "69÷5=" "18÷8="

# Row 5 (table row 17)
Set-Problem 17 1 "95÷8=" "98÷8="
Set-Problem 17 2 "84÷7=" "83÷7="
Set-Problem 17 3 "46÷8=" "80÷3="
Set-Problem 17 4 "40÷9=" "50÷3="
Set-Problem 17 5 "34÷5=" "58÷6="
